$wb = $excel.ActiveWorkbook

# --- Sheet "dic_entidades_nomeadas": clear the "Palavra composta" values in A2:A10 ---
$wsEntidades = $wb.Worksheets.Item("dic_entidades_nomeadas")
$wsEntidades.Range("A2:A10").ClearContents()

# --- Sheet "textos_selecionados": remove rows 3:10, keeping only the header + first data row ---
$wsTextos = $wb.Worksheets.Item("textos_selecionados")
$wsTextos.Rows("3:10").Delete()

# --- Restore / update the selections on each sheet ---
$wsEntidades.Range("J15").Select()

$wsSiglas = $wb.Worksheets.Item("dic_siglas")
$wsSiglas.Range("F8").Select()

$wsTextos.Range("G18").Select()

# --- Make "textos_selecionados" the active (selected) tab ---
$wsTextos.Activate()
